$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update existing scalar odds in rows 4, 6, 7, 8 ---
$ws.Range("G4").Value = 1.85
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48

$ws.Range("O6").Value = 1.24
$ws.Range("P6").Value = 3.85
$ws.Range("S6").Value = 1.32
$ws.Range("T6").Value = 3.25
$ws.Range("U6").Value = 1.78
$ws.Range("V6").Value = 1.98

$ws.Range("O7").Value = 1.18
$ws.Range("P7").Value = 4.51
$ws.Range("S7").Value = 1.32
$ws.Range("T7").Value = 3.25
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.91

$ws.Range("M8").Value = 1.04
$ws.Range("O8").Value = 1.25

# --- Step 2: insert a new row at position 9 (shifts old rows 9,10 down to 10,11) ---
$ws.Rows.Item(9).Insert()

# --- Step 3: populate new row 9 (WATRNsoI / Al Kholood vs Al Nassr) ---
$ws.Range("A9").Value = "WATRNsoI"
$ws.Range("B9").Value = "25/10/2024"
$ws.Range("C9").Value = "12:05"
$ws.Range("D9").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E9").Value = "Al Kholood"
$ws.Range("F9").Value = "Al Nassr"
$ws.Range("G9").Value = 9
$ws.Range("H9").Value = 5.25
$ws.Range("I9").Value = 1.3
$ws.Range("J9").Value = 7
$ws.Range("K9").Value = 2.63
$ws.Range("L9").Value = 1.73
$ws.Range("M9").Value = 1.02
$ws.Range("N9").Value = 11
$ws.Range("O9").Value = 1.14
$ws.Range("P9").Value = 5
$ws.Range("Q9").Value = 1.5
$ws.Range("R9").Value = 2.5
$ws.Range("S9").Value = 1.25
$ws.Range("T9").Value = 3.75
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.83
$ws.Range("W9").Value = 26
$ws.Range("X9").Value = 41
$ws.Range("Y9").Value = 26
$ws.Range("Z9").Value = 101
$ws.Range("AA9").Value = 51
$ws.Range("AB9").Value = 51
$ws.Range("AC9").Value = 17
$ws.Range("AD9").Value = 10
$ws.Range("AE9").Value = 19
$ws.Range("AF9").Value = 51
$ws.Range("AG9").Value = 151
$ws.Range("AH9").Value = 9.5
$ws.Range("AI9").Value = 7.5
$ws.Range("AJ9").Value = 9.5
$ws.Range("AK9").Value = 8.5
$ws.Range("AL9").Value = 11
$ws.Range("AM9").Value = 23
$ws.Range("AN9").Value = 9.5
$ws.Range("AO9").Value = 41
$ws.Range("AP9").Value = 41
$ws.Range("AQ9").Value = 151
$ws.Range("AR9").Value = 151
$ws.Range("AS9").Value = 500
$ws.Range("AT9").Value = 3.75
$ws.Range("AU9").Value = 9
$ws.Range("AV9").Value = 51
$ws.Range("AW9").Value = 81
$ws.Range("AX9").Value = 3.5
$ws.Range("AY9").Value = 6
$ws.Range("AZ9").Value = 15
$ws.Range("BA9").Value = 15
$ws.Range("BB9").Value = 34
$ws.Range("BC9").Value = 101
$ws.Range("BD9").Value = 81

# --- Step 4: add new row 12 (hnBPz8dm / Polissya Zhytomyr vs Kolos Kovalivka) ---
$ws.Range("A12").Value = "hnBPz8dm"
$ws.Range("B12").Value = "25/10/2024"
$ws.Range("C12").Value = "12:00"
$ws.Range("D12").Value = "UKRAINE - PREMIER LEAGUE"
$ws.Range("E12").Value = "Polissya Zhytomyr"
$ws.Range("F12").Value = "Kolos Kovalivka"
$ws.Range("G12").Value = 1.75
$ws.Range("H12").Value = 2.95
$ws.Range("I12").Value = 5.5
$ws.Range("J12").Value = 2.42
$ws.Range("K12").Value = 1.85
$ws.Range("L12").Value = 5.9
$ws.Range("M12").Value = 1.14
$ws.Range("N12").Value = 4.35
$ws.Range("O12").Value = 1.6
$ws.Range("P12").Value = 2.07
$ws.Range("Q12").Value = 2.7
$ws.Range("R12").Value = 1.35
$ws.Range("S12").Value = 1.6
$ws.Range("T12").Value = 2.07
$ws.Range("U12").Value = 2.42
$ws.Range("V12").Value = 1.44
$ws.Range("W12").Value = 4.4
$ws.Range("X12").Value = 6.5
$ws.Range("Y12").Value = 9.5
$ws.Range("Z12").Value = 13.5
$ws.Range("AA12").Value = 20
$ws.Range("AB12").Value = 55
$ws.Range("AC12").Value = 4.65
$ws.Range("AD12").Value = 6.5
$ws.Range("AE12").Value = 27
$ws.Range("AF12").Value = 250
$ws.Range("AG12").Value = 67
$ws.Range("AH12").Value = 9.25
$ws.Range("AI12").Value = 29
$ws.Range("AJ12").Value = 21
$ws.Range("AK12").Value = 120
$ws.Range("AL12").Value = 90
$ws.Range("AM12").Value = 120
$ws.Range("AN12").Value = 3.25
$ws.Range("AO12").Value = 9
$ws.Range("AP12").Value = 26
$ws.Range("AQ12").Value = 37
$ws.Range("AR12").Value = 100
$ws.Range("AS12").Value = 500
$ws.Range("AT12").Value = 2.07
$ws.Range("AU12").Value = 9.25
$ws.Range("AV12").Value = 120
$ws.Range("AW12").Value = 81
$ws.Range("AX12").Value = 6.6
$ws.Range("AY12").Value = 37
$ws.Range("AZ12").Value = 55
$ws.Range("BA12").Value = 300
$ws.Range("BB12").Value = 450
$ws.Range("BC12").Value = 67
$ws.Range("BD12").Value = 81
